# feat: Add SFX of Ship and Bullet resource file
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SFX")

# Row 7 - AllyShipShooting
$ws.Range("A7").Value = "Player&EnemyShipVariety"
$ws.Range("B7").Value = "Bullet"
$ws.Range("C7").Value = "AllyShipShooting"
$ws.Range("D7").Value = "LaserShootAlly1.wav"
$ws.Range("E7").Value = "M"
$ws.Range("F7").Value = "O"
$ws.Range("G7").Formula = '=_xlfn.CONCAT("SFX_",A7,"_",B7,"_",C7,"_",D7)'

# Row 8 - EnemyShipDestroy
$ws.Range("A8").Value = "Player&EnemyShipVariety"
$ws.Range("B8").Value = "EnemyShip"
$ws.Range("C8").Value = "EnemyShipDestroy"
$ws.Range("D8").Value = "Deleted.wav"
$ws.Range("E8").Value = "S"
$ws.Range("F8").Value = "O"
$ws.Range("G8").Formula = '=_xlfn.CONCAT("SFX_",A8,"_",B8,"_",C8,"_",D8)'

# Row 9 - EnemyShipShooting
$ws.Range("A9").Value = "Player&EnemyShipVariety"
$ws.Range("B9").Value = "EnemyShip"
$ws.Range("C9").Value = "EnemyShipShooting"
$ws.Range("D9").Value = "LaserShootEnemy1.wav"
$ws.Range("E9").Value = "M"
$ws.Range("F9").Value = "O"
$ws.Range("G9").Formula = '=_xlfn.CONCAT("SFX_",A9,"_",B9,"_",C9,"_",D9)'

# Row 10 - EnemyShipSpecial
$ws.Range("A10").Value = "Player&EnemyShipVariety"
$ws.Range("B10").Value = "EnemyShip"
$ws.Range("C10").Value = "EnemyShipSpecial"
$ws.Range("D10").Value = "ShipAppear.wav"
$ws.Range("E10").Value = "S"
$ws.Range("F10").Value = "O"
$ws.Range("G10").Formula = '=_xlfn.CONCAT("SFX_",A10,"_",B10,"_",C10,"_",D10)'

# Row 11 - EnemyShipSpecialDestroy
$ws.Range("A11").Value = "Player&EnemyShipVariety"
$ws.Range("B11").Value = "EnemyShip"
$ws.Range("C11").Value = "EnemyShipSpecialDestroy"
$ws.Range("D11").Value = "ExplosionShort.wav"
$ws.Range("E11").Value = "S"
$ws.Range("F11").Value = "O"
$ws.Range("G11").Formula = '=_xlfn.CONCAT("SFX_",A11,"_",B11,"_",C11,"_",D11)'

# Row 12 - AllyShipDestroy
$ws.Range("A12").Value = "Player&EnemyShipVariety"
$ws.Range("B12").Value = "Ship"
$ws.Range("C12").Value = "AllyShipDestroy"
$ws.Range("D12").Value = "ExplosionLong.wav"
$ws.Range("E12").Value = "S"
$ws.Range("F12").Value = "O"
$ws.Range("G12").Formula = '=_xlfn.CONCAT("SFX_",A12,"_",B12,"_",C12,"_",D12)'
$ws.Range("C12").Font.Color = 0

# Row 13 - AllyShipDestroyLives
$ws.Range("A13").Value = "Player&EnemyShipVariety"
$ws.Range("B13").Value = "Ship"
$ws.Range("C13").Value = "AllyShipDestroyLives"
$ws.Range("D13").Value = "DescendingScales1.wav"
$ws.Range("E13").Value = "S"
$ws.Range("F13").Value = "O"
$ws.Range("G13").Formula = '=_xlfn.CONCAT("SFX_",A13,"_",B13,"_",C13,"_",D13)'

# Remove the old leftover placeholder row
$ws.Range("A17:G17").EntireRow.Delete()
